# Updates the "cryptos" price list: refreshed Price (D) and Volume(1h) (E)
# values for most rows, plus a swap of the RenderToken-adjacent entries so
# row 31 becomes Binance-PegBSC-USD and row 32 becomes Fetch.AI (with their
# own refreshed Price/Volume), matching the upstream GitHub Actions refresh.
#
# NumberFormat is forced to Text ("@") immediately before each Price/Volume
# Value assignment: these columns store plain text (e.g. "61.625.95",
# "1.00", "  +1.37%  ") and, without this, Excel's COM layer would parse
# numeric-looking strings (like "1.00" or "7.80") into actual numbers and
# silently drop the trailing zero / formatting.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = "@"
$ws.Range('D2').Value = '61.625.95'
$ws.Range('E2').NumberFormat = "@"
$ws.Range('E2').Value = '  +1.37%  '
$ws.Range('D3').NumberFormat = "@"
$ws.Range('D3').Value = '3.447.24'
$ws.Range('E3').NumberFormat = "@"
$ws.Range('E3').Value = '  +2.25%  '
$ws.Range('D4').NumberFormat = "@"
$ws.Range('D4').Value = '1.00'
$ws.Range('E4').NumberFormat = "@"
$ws.Range('E4').Value = '  -0.05%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '579.98'
$ws.Range('E5').NumberFormat = "@"
$ws.Range('E5').Value = '  +1.40%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '149.39'
$ws.Range('E6').NumberFormat = "@"
$ws.Range('E6').Value = '  +9.16%  '
$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '3.449.01'
$ws.Range('E7').NumberFormat = "@"
$ws.Range('E7').Value = '  +2.36%  '
$ws.Range('E8').NumberFormat = "@"
$ws.Range('E8').Value = '  +0.05%  '
$ws.Range('E9').NumberFormat = "@"
$ws.Range('E9').Value = '  +1.07%  '
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '7.80'
$ws.Range('E10').NumberFormat = "@"
$ws.Range('E10').Value = '  +1.86%  '
$ws.Range('E11').NumberFormat = "@"
$ws.Range('E11').Value = '  +3.42%  '
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '0.391'
$ws.Range('E12').NumberFormat = "@"
$ws.Range('E12').Value = '  +1.66%  '
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '4.036.47'
$ws.Range('E13').NumberFormat = "@"
$ws.Range('E13').Value = '  +2.26%  '
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '27.89'
$ws.Range('E14').NumberFormat = "@"
$ws.Range('E14').Value = '  +7.47%  '
$ws.Range('E15').NumberFormat = "@"
$ws.Range('E15').Value = '  -0.45%  '
$ws.Range('E16').NumberFormat = "@"
$ws.Range('E16').Value = '  +2.30%  '
$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '3.455.05'
$ws.Range('E17').NumberFormat = "@"
$ws.Range('E17').Value = '  +2.52%  '
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '61.724.87'
$ws.Range('E18').NumberFormat = "@"
$ws.Range('E18').Value = '  +1.25%  '
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '6.28'
$ws.Range('E19').NumberFormat = "@"
$ws.Range('E19').Value = '  +8.51%  '
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '14.36'
$ws.Range('E20').NumberFormat = "@"
$ws.Range('E20').Value = '  +3.05%  '
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '9.49'
$ws.Range('E21').NumberFormat = "@"
$ws.Range('E21').Value = '  +0.85%  '
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '389.54'
$ws.Range('E22').NumberFormat = "@"
$ws.Range('E22').Value = '  +4.21%  '
$ws.Range('E23').NumberFormat = "@"
$ws.Range('E23').Value = '  +2.78%  '
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '3.589.15'
$ws.Range('E24').NumberFormat = "@"
$ws.Range('E24').Value = '  +2.02%  '
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '73.17'
$ws.Range('E25').NumberFormat = "@"
$ws.Range('E25').Value = '  +3.00%  '
$ws.Range('E26').NumberFormat = "@"
$ws.Range('E26').Value = '  +0.14%  '
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '5.76'
$ws.Range('E27').NumberFormat = "@"
$ws.Range('E27').Value = '  +0.36%  '
$ws.Range('E28').NumberFormat = "@"
$ws.Range('E28').Value = '  +0.26%  '
$ws.Range('E29').NumberFormat = "@"
$ws.Range('E29').Value = '  +2.46%  '
$ws.Range('E30').NumberFormat = "@"
$ws.Range('E30').Value = '  +3.60%  '
$ws.Range('B31').Value = 'Binance-PegBSC-USD'
$ws.Range('C31').Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '1.00'
$ws.Range('E31').NumberFormat = "@"
$ws.Range('E31').Value = '  +0.12%  '
$ws.Range('B32').Value = 'Fetch.AI'
$ws.Range('C32').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '1.54'
$ws.Range('E32').NumberFormat = "@"
$ws.Range('E32').Value = '  -13.42%  '
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '8.24'
$ws.Range('E33').NumberFormat = "@"
$ws.Range('E33').Value = '  +1.68%  '
$ws.Range('E34').NumberFormat = "@"
$ws.Range('E34').Value = '  +1.36%  '
$ws.Range('E35').NumberFormat = "@"
$ws.Range('E35').Value = '  -0.06%  '
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '24.04'
$ws.Range('E36').NumberFormat = "@"
$ws.Range('E36').Value = '  +1.67%  '
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '7.02'
$ws.Range('E37').NumberFormat = "@"
$ws.Range('E37').Value = '  +2.37%  '
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '5.21'
$ws.Range('E38').NumberFormat = "@"
$ws.Range('E38').Value = '  +0.58%  '
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '1.56'
$ws.Range('E39').NumberFormat = "@"
$ws.Range('E39').Value = '  +1.54%  '
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '166.35'
$ws.Range('E40').NumberFormat = "@"
$ws.Range('E40').Value = '  +1.11%  '
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '0.0785'
$ws.Range('E41').NumberFormat = "@"
$ws.Range('E41').Value = '  +3.26%  '
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '26.83'
$ws.Range('E42').NumberFormat = "@"
$ws.Range('E42').Value = '  +11.20%  '
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '0.792'
$ws.Range('E43').NumberFormat = "@"
$ws.Range('E43').Value = '  +2.18%  '
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '4.51'
$ws.Range('E44').NumberFormat = "@"
$ws.Range('E44').Value = '  +2.66%  '
$ws.Range('E45').NumberFormat = "@"
$ws.Range('E45').Value = '  -0.03%  '
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '42.32'
$ws.Range('E46').NumberFormat = "@"
$ws.Range('E46').Value = '  +2.02%  '
$ws.Range('E47').NumberFormat = "@"
$ws.Range('E47').Value = '  +1.18%  '
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '2.610.32'
$ws.Range('E48').NumberFormat = "@"
$ws.Range('E48').Value = '  +6.36%  '
$ws.Range('E49').NumberFormat = "@"
$ws.Range('E49').Value = '  -1.80%  '
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '6.95'
$ws.Range('E50').NumberFormat = "@"
$ws.Range('E50').Value = '  +2.40%  '
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '23.21'
$ws.Range('E51').NumberFormat = "@"
$ws.Range('E51').Value = '  +1.40%  '
